# Generate Report for handback
# Marks the a.md.md / b.md.md rows as handed back (in sync with en-US),
# filling in the "Latest Target File" / "Latest Handback File" columns
# and stamping a "Latest Handback DateTime", on both the zh-cn and de-de
# language sheets plus the Overview roll-up sheet.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the status text changes (it mirrors the per-
# language Status column for a.md.md / b.md.md).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcnMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/1dd947e83bbc9167b458b976f699f0bcaf2e56a5/e2e/a.md.md"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8db822c8f780c1eda7dd910c66b34e76585367d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
$zhcnXlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
$zhcnHandbackDatetime = "2016-01-25 03:27:26"

# Row 2 (a.md.md)
$zhcn.Range("B2").Value = $statusHandedBack
$zhcn.Range("E2").Value = "a.md.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), $zhcnMdUrl, "", "", "a.md.md")
$zhcn.Range("E2").Font.Underline = $True
$zhcn.Range("E2").Font.Color = 15570276
$zhcn.Range("F2").Value = $zhcnXlfName
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhcnXlfUrl, "", "", $zhcnXlfName)
$zhcn.Range("F2").Font.Underline = $True
$zhcn.Range("F2").Font.Color = 15570276
$zhcn.Range("G2").Value = $zhcnHandbackDatetime

# Row 3 (b.md.md) -- Target/Handback columns mirror the a.md.md ones,
# matching the source report.
$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("E3").Value = "a.md.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), $zhcnMdUrl, "", "", "a.md.md")
$zhcn.Range("E3").Font.Underline = $True
$zhcn.Range("E3").Font.Color = 15570276
$zhcn.Range("F3").Value = $zhcnXlfName
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhcnXlfUrl, "", "", $zhcnXlfName)
$zhcn.Range("F3").Font.Underline = $True
$zhcn.Range("F3").Font.Color = 15570276
$zhcn.Range("G3").Value = $zhcnHandbackDatetime

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dedeMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/1dd947e83bbc9167b458b976f699f0bcaf2e56a5/e2e/a.md.md"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e23a9f3155da1bbb3e9aca31a974a0e3637d3d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
$dedeXlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
$dedeHandbackDatetime = "2016-01-25 03:27:45"

# Row 2 (a.md.md)
$dede.Range("B2").Value = $statusHandedBack
$dede.Range("E2").Value = "a.md.md"
$dede.Hyperlinks.Add($dede.Range("E2"), $dedeMdUrl, "", "", "a.md.md")
$dede.Range("E2").Font.Underline = $True
$dede.Range("E2").Font.Color = 15570276
$dede.Range("F2").Value = $dedeXlfName
$dede.Hyperlinks.Add($dede.Range("F2"), $dedeXlfUrl, "", "", $dedeXlfName)
$dede.Range("F2").Font.Underline = $True
$dede.Range("F2").Font.Color = 15570276
$dede.Range("G2").Value = $dedeHandbackDatetime

# Row 3 (b.md.md)
$dede.Range("B3").Value = $statusHandedBack
$dede.Range("E3").Value = "a.md.md"
$dede.Hyperlinks.Add($dede.Range("E3"), $dedeMdUrl, "", "", "a.md.md")
$dede.Range("E3").Font.Underline = $True
$dede.Range("E3").Font.Color = 15570276
$dede.Range("F3").Value = $dedeXlfName
$dede.Hyperlinks.Add($dede.Range("F3"), $dedeXlfUrl, "", "", $dedeXlfName)
$dede.Range("F3").Font.Underline = $True
$dede.Range("F3").Font.Color = 15570276
$dede.Range("G3").Value = $dedeHandbackDatetime
